$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "테디노트"

$ws.Range("D26").Value = "인공지능 음성 생성 연구: 음성 분류 솔루션"

$ws.Range("D36").Value = "dmqm_seminar"

$ws.Range("D37").Value = "dsba_seminar"

$ws.Range("D42").Value = "tensorflow CUBLAS_STATUS_ALLOC_FAILED 오류 대처 및 메모리 할당"
$ws.Range("E42").Value = "https://kjk92.tistory.com/87"

$ws.Range("D49").Value = "taeu"

$ws.Range("D50").Value = "바죠"

$ws.Range("D51").Value = "코딩이 잘 안 될 때는 잠시 쉬자"
$ws.Range("E51").Value = "https://bskyvision.com/entry/%ED%94%84%EB%A1%9C%EA%B7%B8%EB%9E%98%EB%B0%8D%EC%9D%B4-%EC%9E%98-%EC%95%88-%ED%92%80%EB%A6%B4-%EB%95%8C%EB%8A%94-%EC%9E%A0%EC%8B%9C-%EC%89%AC%EC%9E%90"
